# Updated cryptos list on Fri Nov  1 03:54:23 UTC 2024 with GitHub Actions
# Refresh Price (column D) and Volume(1h) (column E) figures for each coin row,
# and reorder/replace a few rows: ARBITRUM <-> Filecoin swap (rows 48/49),
# and Mantle -> Cronos (row 51).
#
# Numeric-looking Price strings are forced to stay text (NumberFormat "@")
# so Excel does not silently convert them to floating point numbers -
# matching the workbook convention of storing Price as a string.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.395.50"
$ws.Range("E2").Value = "  -4.13%  "
$ws.Range("D3").Value = "2.500.86"
$ws.Range("E3").Value = "  -5.93%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.16%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.516"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.07%  "
$ws.Range("D9").Value = "2.500.63"
$ws.Range("E9").Value = "  -5.91%  "
$ws.Range("E10").Value = "  -6.88%  "
$ws.Range("E11").Value = "  -0.59%  "
$ws.Range("E12").Value = "  -5.10%  "
$ws.Range("E13").Value = "  -2.65%  "
$ws.Range("D14").Value = "2.961.30"
$ws.Range("E14").Value = "  -5.89%  "
$ws.Range("D15").Value = "69.417.54"
$ws.Range("E15").Value = "  -3.94%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000175"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.51%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.87"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.65%  "
$ws.Range("D18").Value = "2.503.66"
$ws.Range("E18").Value = "  -5.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.37"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -8.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.73"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "351.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.94"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.76%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.98"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.93%  "
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("E25").Value = "  -3.52%  "
$ws.Range("E26").Value = "  -5.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.996"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.40%  "
$ws.Range("E30").Value = "  -5.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.89"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "483.16"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.14%  "
$ws.Range("E33").Value = "  +1.08%  "
$ws.Range("E34").Value = "  -2.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("E36").Value = "  -2.36%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "152.02"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.88"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.56"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.43%  "
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("E41").Value = "  -2.27%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.319"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.62"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.04%  "
$ws.Range("E44").Value = "  -13.81%  "
$ws.Range("E45").Value = "  -8.25%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "38.22"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.46%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "143.45"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.52%  "
$ws.Range("B48").Value = "Filecoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.54"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.14%  "
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.529"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.94%  "
$ws.Range("E50").Value = "  -5.37%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0730"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.49%  "
